# Region_Estimaiton_Result.xlsx - "add data of AHP"
#
# Both sheets ("pre" and "standard") held a 4x3 table (row labels
# Y / Y_Import / r_c / r_p down column A, column headers China / EU / India
# across row 1). The table is transposed in place: the new row labels
# (down column A) become China / EU / India, the new column headers
# (across row 1) become Y / Y_Import / r_c / r_p, and a new "Region"
# label is added at A1. The data values move along with their row/column
# labels, so the sheet shrinks from 5 rows to 4 rows (A1:D5 -> A1:E4).

$wb = $excel.ActiveWorkbook

function Set-TransposedSheet {
    param($ws)

    # Row 5 no longer exists after the transpose (4 labels instead of 5) -
    # delete it outright so the sheet dimension shrinks to A1:E4 and the
    # remaining rows keep their existing formatting untouched.
    $ws.Rows.Item(5).Delete()

    # A1 and the new E1 header need the same bold/bordered/centered style
    # already used by the other header cells - copy it over before writing
    # the literal values.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $ws.Range("E1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Header row: Region | Y | Y_Import | r_c | r_p
    $ws.Range("A1").Value = "Region"
    $ws.Range("B1").Value = "Y"
    $ws.Range("C1").Value = "Y_Import"
    $ws.Range("D1").Value = "r_c"
    $ws.Range("E1").Value = "r_p"

    # Row labels: China / EU / India
    $ws.Range("A2").Value = "China"
    $ws.Range("A3").Value = "EU"
    $ws.Range("A4").Value = "India"
}

$wsPre = $wb.Worksheets.Item("pre")
Set-TransposedSheet $wsPre

$wsPre.Range("B2").Value = -9.932487581726193
$wsPre.Range("C2").Value = 6.97521704121815
$wsPre.Range("D2").Value = -5.032224880158286
$wsPre.Range("E2").Value = -3.436828665498592

$wsPre.Range("B3").Value = 24.98504778026656
$wsPre.Range("C3").Value = -14.14170446634398
$wsPre.Range("D3").Value = 10.84554644021556
$wsPre.Range("E3").Value = 8.7831677131497

$wsPre.Range("B4").Value = -15.05256019854037
$wsPre.Range("C4").Value = 7.166487425125831
$wsPre.Range("D4").Value = -5.813321560057277
$wsPre.Range("E4").Value = -5.346339047651106

$wsStandard = $wb.Worksheets.Item("standard")
Set-TransposedSheet $wsStandard

$wsStandard.Range("B2").Value = -0.5583085299284382
$wsStandard.Range("C2").Value = 0.697521704121815
$wsStandard.Range("D2").Value = -0.6556144016474311
$wsStandard.Range("E2").Value = -0.549069353459124

$wsStandard.Range("B3").Value = 1.404418095830945
$wsStandard.Range("C3").Value = -1.414170446634398
$wsStandard.Range("D3").Value = 1.412992584647309
$wsStandard.Range("E3").Value = 1.403202977790146

$wsStandard.Range("B4").Value = -0.846109565902507
$wsStandard.Range("C4").Value = 0.7166487425125831
$wsStandard.Range("D4").Value = -0.7573781829998781
$wsStandard.Range("E4").Value = -0.8541336243310218
